$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (Q0) ---
$ws.Range("B2").Value = 0.1805204829727294
$ws.Range("C2").Value = 0.9134701371387803
$ws.Range("D2").Value = 2.313862988505206
$ws.Range("E2").Value = 1.521138714419302
$ws.Range("F2").Value = 1.567404881478029

# --- Row 3 (Q1) ---
$ws.Range("B3").Value = 0.1654443250838863
$ws.Range("C3").Value = 0.8321534822147496
$ws.Range("D3").Value = 2.430194309093069
$ws.Range("E3").Value = 1.55890805023679
$ws.Range("F3").Value = 1.613399420919432

# --- Row 4 (Q2) ---
$ws.Range("B4").Value = 0.3342101850538108
$ws.Range("C4").Value = 0.8990974440469651
$ws.Range("D4").Value = 3.047122944578856
$ws.Range("E4").Value = 1.745601026746621
$ws.Range("F4").Value = 1.789492512149283

# --- Row 5 (Q3) ---
$ws.Range("B5").Value = 0.349057914304527
$ws.Range("C5").Value = 0.9234061813053621
$ws.Range("D5").Value = 2.965460108284523
$ws.Range("E5").Value = 1.722051134050474
$ws.Range("F5").Value = 1.768609778560688

# --- Row 6 (Q4) ---
$ws.Range("B6").Value = 0.4342163702162936
$ws.Range("C6").Value = 1.060274838210373
$ws.Range("D6").Value = 3.344243897922033
$ws.Range("E6").Value = 1.828727398471963
$ws.Range("F6").Value = 1.872520595275618

# --- Row 7 (Q5) ---
$ws.Range("B7").Value = -0.06951841543571383
$ws.Range("C7").Value = 0.5087111458309537
$ws.Range("D7").Value = 0.3139214969299044
$ws.Range("E7").Value = 0.560286977298156
$ws.Range("F7").Value = 0.5896819250247234
$ws.Range("G7").Value = 9

# --- Row 8 (Q6) ---
$ws.Range("B8").Value = -0.04657217614193387
$ws.Range("C8").Value = 0.7130288563221177
$ws.Range("D8").Value = 0.6707611900940872
$ws.Range("E8").Value = 0.8190001160525481
$ws.Range("F8").Value = 0.8957179617514601
$ws.Range("G8").Value = 6

# --- Row 9 (Q7) ---
$ws.Range("B9").Value = -0.3788153344042176
$ws.Range("C9").Value = 0.3788153344042176
$ws.Range("D9").Value = 0.2155603798649167
$ws.Range("E9").Value = 0.4642848046887995
$ws.Range("F9").Value = 0.3287688906020555
$ws.Range("G9").Value = 3

# --- Row 10 (new: Q8) ---
$ws.Range("A10").Value = "Q8"
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$ws.Range("B10").Value = 0.410128014204804
$ws.Range("C10").Value = 0.410128014204804
$ws.Range("D10").Value = 0.168204988035576
$ws.Range("E10").Value = 0.410128014204804
$ws.Range("G10").Value = 1
